$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Data entry: "gp" (Bewertung/Rating, rows 16-29) and "Notizen" (Anzahl
# Bewertungen, rows 16-29) columns that were still blank for the "With the
# Beatles" album block get filled in. Row 20 and row 26 only received a
# ratings-count (column C); column B stayed empty for those two rows.
$data = @{
    16 = @{ B = 5;   C = 8  }
    17 = @{ B = 4.5; C = 3  }
    18 = @{ B = 5;   C = 30 }
    19 = @{ B = 4.5; C = 4  }
    20 = @{ C = 0 }
    21 = @{ B = 5;   C = 13 }
    22 = @{ B = 5;   C = 9  }
    23 = @{ B = 5;   C = 9  }
    24 = @{ B = 5;   C = 8  }
    25 = @{ B = 4.5; C = 5  }
    26 = @{ C = 0 }
    27 = @{ B = 4.5; C = 3  }
    28 = @{ B = 5;   C = 5  }
    29 = @{ B = 5;   C = 4  }
}

foreach ($row in ($data.Keys | Sort-Object)) {
    $vals = $data[$row]
    if ($vals.ContainsKey("B")) {
        $ws.Cells.Item($row, 2).Value = $vals["B"]
    }
    if ($vals.ContainsKey("C")) {
        $ws.Cells.Item($row, 3).Value = $vals["C"]
    }
}

# Highlight the "I Wanna Be Your Man" title (A26) with a yellow fill - a
# little "notiz" marker on that row.
$ws.Range("A26").Interior.Color = 65535

# --- View state: restore the scroll/zoom/selection the workbook was left
# in (zoomed in around row 10, cell C29 selected).
$ws.Range("C29").Select()
$win = $excel.ActiveWindow
$win.Zoom = 122
$win.ScrollRow = 10
$win.ScrollColumn = 1

# Window geometry of the Excel app itself when the file was last saved.
$win.Left = 2440
$win.Top = 2500
$win.Width = 23740
$win.Height = 12640
